$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two rows holding the standalone "5840963 - Daniela Camargo Vernilli" and
# "5840820 - Gustavo Aristides Santana Martinez" entries (old rows 13 & 14) are
# removed; everything below shifts up by two rows (dimension becomes A1:C24).
$ws.Rows.Item(13).EntireRow.Delete()
$ws.Rows.Item(13).EntireRow.Delete()

# Re-populate the cells whose content was rearranged as part of this edit.
$ws.Range("B10").Value = "5840963 - Daniela Camargo Vernilli"
$ws.Range("C10").Value = "5840963 - Daniela Camargo Vernilli"

$ws.Range("B13").Value = "01/01/2017"
$ws.Range("C13").Value = "01/01/2017"

$ws.Range("B15").Value = "5840963 - Daniela Camargo Vernilli"
$ws.Range("C15").Value = "5840963 - Daniela Camargo Vernilli"

$ws.Range("B18").Value = "5840820 - Gustavo Aristides Santana Martinez"
$ws.Range("C18").Value = "5840820 - Gustavo Aristides Santana Martinez"

$ws.Range("B19").Value = "Duas provas"
$ws.Range("C19").Value = "Duas provas"

$ws.Range("B20").Value = "Serão aplicadas duas avaliações (P1 e P2) que comporão a nota final (NF). A nota final será calculada através da expressão: NF= (P1+P2)/2"
$ws.Range("C20").Value = "Serão aplicadas duas avaliações (P1 e P2) que comporão a nota final (NF). A nota final será calculada através da expressão: NF= (P1+P2)/2"

$ws.Range("B21").Value = "Para o aluno que obtiver Nota Final maior ou igual a 3,0 e menor do que 5,0, será aplicada uma avaliação de recuperação (NR), com pontuação de 0 a 10, que levará ao cálculo da média final(MF) através da seguinte expressão: MF=(NF+NR)/2.onde: NF=Nota Final e NR=Nota da Prova de Recuperação."
$ws.Range("C21").Value = "Para o aluno que obtiver Nota Final maior ou igual a 3,0 e menor do que 5,0, será aplicada uma avaliação de recuperação (NR), com pontuação de 0 a 10, que levará ao cálculo da média final(MF) através da seguinte expressão: MF=(NF+NR)/2.onde: NF=Nota Final e NR=Nota da Prova de Recuperação."
